# Apply updated cryptocurrency price/volume figures to the "cryptos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) whose new text looks like a plain decimal number must be
# written with a leading apostrophe so Excel keeps them as text (e.g. '1.00', '42.90')
# instead of silently converting them to a number and losing the trailing zero / exact
# decimal representation that the source data uses.

# Row 2
$ws.Range('D2').Value = '64.961.49'
$ws.Range('E2').Value = '  +0.13%  '

# Row 3
$ws.Range('D3').Value = '3.520.87'
$ws.Range('E3').Value = '  -0.85%  '

# Row 4
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').Value = '''592.39'
$ws.Range('E5').Value = '  -1.11%  '

# Row 6
$ws.Range('D6').Value = '''134.05'
$ws.Range('E6').Value = '  -1.60%  '

# Row 7
$ws.Range('D7').Value = '3.520.71'

# Row 8
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').Value = '''0.491'
$ws.Range('E9').Value = '  -1.02%  '

# Row 10
$ws.Range('D10').Value = '''0.124'
$ws.Range('E10').Value = '  +1.00%  '

# Row 11
$ws.Range('D11').Value = '''7.14'
$ws.Range('E11').Value = '  +2.95%  '

# Row 12
$ws.Range('D12').Value = '''0.387'
$ws.Range('E12').Value = '  +0.29%  '

# Row 13
$ws.Range('D13').Value = '4.125.63'
$ws.Range('E13').Value = '  -0.76%  '

# Row 14
$ws.Range('D14').Value = '''27.61'
$ws.Range('E14').Value = '  +2.01%  '

# Row 15
$ws.Range('E15').Value = '  -0.22%  '

# Row 16
$ws.Range('E16').Value = '  +0.45%  '

# Row 17
$ws.Range('D17').Value = '3.521.68'
$ws.Range('E17').Value = '  -1.14%  '

# Row 18
$ws.Range('D18').Value = '65.008.33'
$ws.Range('E18').Value = '  +0.32%  '

# Row 19
$ws.Range('D19').Value = '''10.13'
$ws.Range('E19').Value = '  +1.02%  '

# Row 20
$ws.Range('D20').Value = '''14.45'
$ws.Range('E20').Value = '  +0.33%  '

# Row 21
$ws.Range('D21').Value = '''5.71'
$ws.Range('E21').Value = '  -1.90%  '

# Row 22
$ws.Range('D22').Value = '''391.83'
$ws.Range('E22').Value = '  +0.62%  '

# Row 23
$ws.Range('D23').Value = '''0.579'
$ws.Range('E23').Value = '  +0.73%  '

# Row 24
$ws.Range('D24').Value = '3.665.62'
$ws.Range('E24').Value = '  -0.87%  '

# Row 25
$ws.Range('D25').Value = '''74.74'
$ws.Range('E25').Value = '  +0.79%  '

# Row 26
$ws.Range('E26').Value = '  +0.52%  '

# Row 27
$ws.Range('E27').Value = '  -3.50%  '

# Row 28
$ws.Range('D28').Value = '''7.72'
$ws.Range('E28').Value = '  +0.90%  '

# Row 29
$ws.Range('E29').Value = '  +10.82%  '

# Row 30
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  -0.01%  '

# Row 31
$ws.Range('D31').Value = '''2.27'
$ws.Range('E31').Value = '  -0.70%  '

# Row 32
$ws.Range('D32').Value = '''8.36'
$ws.Range('E32').Value = '  +0.29%  '

# Row 33
$ws.Range('D33').Value = '3.525.92'
$ws.Range('E33').Value = '  -1.01%  '

# Row 34
$ws.Range('D34').Value = '''24.10'
$ws.Range('E34').Value = '  +0.47%  '

# Row 35
$ws.Range('E35').Value = '  -0.01%  '

# Row 36
$ws.Range('D36').Value = '''0.144'
$ws.Range('E36').Value = '  -0.68%  '

# Row 37
$ws.Range('D37').Value = '''5.32'
$ws.Range('E37').Value = '  +6.65%  '

# Row 38
$ws.Range('E38').Value = '  +3.16%  '

# Row 39
$ws.Range('D39').Value = '''7.01'
$ws.Range('E39').Value = '  +1.11%  '

# Row 40
$ws.Range('D40').Value = '''168.67'
$ws.Range('E40').Value = '  -0.98%  '

# Row 41
$ws.Range('D41').Value = '''0.0815'
$ws.Range('E41').Value = '  +1.21%  '

# Row 42
$ws.Range('D42').Value = '''0.824'
$ws.Range('E42').Value = '  -0.30%  '

# Row 43
$ws.Range('D43').Value = '''1.27'
$ws.Range('E43').Value = '  +5.70%  '

# Row 44
$ws.Range('D44').Value = '''25.94'
$ws.Range('E44').Value = '  -3.05%  '

# Row 45
$ws.Range('D45').Value = '''42.90'
$ws.Range('E45').Value = '  +0.70%  '

# Row 46
$ws.Range('E46').Value = '  +0.01%  '

# Row 47
$ws.Range('D47').Value = '''4.44'
$ws.Range('E47').Value = '  -0.11%  '

# Row 48
$ws.Range('E48').Value = '  +0.66%  '

# Row 49
$ws.Range('D49').Value = '''6.94'
$ws.Range('E49').Value = '  +0.80%  '

# Row 50
$ws.Range('D50').Value = '2.415.26'
$ws.Range('E50').Value = '  -1.52%  '

# Row 51
$ws.Range('D51').Value = '''0.908'
$ws.Range('E51').Value = '  +6.26%  '
